$d = $word.ActiveDocument

# Paragraph 1 is the "Write Up" title. Paragraph 2 is the first of the
# pre-existing trailing blank paragraphs (no explicit style / pPr).
# Insert three blank paragraphs right before it so the new paragraphs
# don't pick up the "Title" style that InsertParagraphAfter on
# paragraph 1 would otherwise copy onto them.
$anchor = $d.Paragraphs(2).Range
$anchor.InsertParagraphBefore()
$anchor.InsertParagraphBefore()
$anchor.InsertParagraphBefore()

function Set-ParagraphXml($paraIndex, [string]$innerParaXml) {
    # Replace the (currently empty) paragraph at $paraIndex with the
    # OOXML supplied in $innerParaXml via Range.InsertXML. This keeps
    # runs distinct (InsertXML does not coalesce adjacent runs the way
    # successive Range.InsertAfter calls / Range.Text assignment can),
    # which is needed for the "berserk" paragraph below.
    $rng = $word.ActiveDocument.Paragraphs($paraIndex).Range
    $rng.Collapse(1)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

$para2xml = '<w:p><w:r><w:t>This week, we will be taking a look at some really weird measurement behavior from CSS, and using a technique called Box-Sizing to fix it.</w:t></w:r></w:p>'
Set-ParagraphXml 2 $para2xml

$para3xml = '<w:p><w:r><w:t xml:space="preserve">So, if you are being driven crazy from things going </w:t></w:r>'
$para3xml = $para3xml + '<w:r><w:t>b</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t>rserk</w:t></w:r>'
$para3xml = $para3xml + '<w:r><w:t xml:space="preserve"> after adding a border, or padding to a box, you may want to join us for our brand-new article this week entitled:</w:t></w:r></w:p>'
Set-ParagraphXml 3 $para3xml

$para4xml = '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Box-Sizing</w:t></w:r></w:p>'
Set-ParagraphXml 4 $para4xml
